$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row with "Apple" (matches shared string index 1 reused) in A5
$ws.Range("A5").Value = "Apple"

# Select the new cell to match the sheetView selection change in the diff
$ws.Range("A5").Select()
